$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "95.569.09"
$ws.Range("E2").Value = "  +2.21%  "

# Row 3
$ws.Range("D3").Value = "3.599.15"
$ws.Range("E3").Value = "  +5.03%  "

# Row 4
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$ws.Range("D5").Value = "'238.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.88%  "

# Row 6
$ws.Range("D6").Value = "'653.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.04%  "

# Row 7
$ws.Range("D7").Value = "'1.47"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.48%  "

# Row 8
$ws.Range("D8").Value = "'0.404"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.92%  "

# Row 9
$ws.Range("E9").Value = "  -0.05%  "

# Row 10
$ws.Range("D10").Value = "'1.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.61%  "

# Row 11
$ws.Range("D11").Value = "3.600.58"
$ws.Range("E11").Value = "  +5.10%  "

# Row 12
$ws.Range("D12").Value = "'42.80"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.62%  "

# Row 13
$ws.Range("D13").Value = "'0.200"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.42%  "

# Row 14
$ws.Range("D14").Value = "'6.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.62%  "

# Row 15
$ws.Range("D15").Value = "4.270.53"
$ws.Range("E15").Value = "  +4.79%  "

# Row 16
$ws.Range("D16").Value = "95.504.56"
$ws.Range("E16").Value = "  +2.39%  "

# Row 17
$ws.Range("D17").Value = "'0.0000255"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.04%  "

# Row 18
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.597.12"
$ws.Range("E18").Value = "  +4.81%  "

# Row 19
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "'7.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.80%  "

# Row 20
$ws.Range("D20").Value = "'12.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.33%  "

# Row 21
$ws.Range("D21").Value = "'17.96"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.76%  "

# Row 22
$ws.Range("D22").Value = "'3.61"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.06%  "

# Row 23
$ws.Range("D23").Value = "'0.488"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.36%  "

# Row 24
$ws.Range("D24").Value = "'509.64"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.34%  "

# Row 25
$ws.Range("D25").Value = "'0.0000196"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.91%  "

# Row 26
$ws.Range("D26").Value = "'6.62"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.32%  "

# Row 27
$ws.Range("D27").Value = "'96.82"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.81%  "

# Row 28
$ws.Range("D28").Value = "'12.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.90%  "

# Row 29
$ws.Range("D29").Value = "3.800.91"
$ws.Range("E29").Value = "  +5.25%  "

# Row 30
$ws.Range("D30").Value = "'3.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +12.44%  "

# Row 31
$ws.Range("D31").Value = "'11.33"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.21%  "

# Row 32
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "'0.140"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.04%  "

# Row 33
$ws.Range("B33").Value = "Dai"
$ws.Range("C33").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D33").Value = "'0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.09%  "

# Row 34
$ws.Range("D34").Value = "'0.997"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.67%  "

# Row 35
$ws.Range("E35").Value = "  +2.92%  "

# Row 36
$ws.Range("D36").Value = "'31.94"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.35%  "

# Row 37
$ws.Range("D37").Value = "'0.561"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.44%  "

# Row 38
$ws.Range("D38").Value = "'8.22"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +9.45%  "

# Row 39
$ws.Range("D39").Value = "'573.66"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.72%  "

# Row 40
$ws.Range("D40").Value = "'1.48"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.02%  "

# Row 41
$ws.Range("E41").Value = "  -0.02%  "

# Row 42
$ws.Range("D42").Value = "'0.150"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.58%  "

# Row 43
$ws.Range("D43").Value = "'0.921"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.84%  "

# Row 44
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'34.62"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +36.56%  "

# Row 45
$ws.Range("B45").Value = "ImmutableX"
$ws.Range("C45").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D45").Value = "'1.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.19%  "

# Row 46
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").Value = "'5.70"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.29%  "

# Row 47
$ws.Range("B47").Value = "WhiteBITCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D47").Value = "'23.77"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.32%  "

# Row 48
$ws.Range("D48").Value = "'2.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.94%  "

# Row 49
$ws.Range("D49").Value = "'0.0415"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.07%  "

# Row 50
$ws.Range("D50").Value = "'3.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.89%  "

# Row 51
$ws.Range("D51").Value = "'53.78"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.23%  "
